# "time card and backlog update"
# Adds Week 5 ("G") and Week 6 ("H") hours for each person on Sheet1,
# adds a per-person TOTAL column (R) with a SUM formula, and moves the
# active-cell selection to reflect the newly entered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New week-5 / week-6 hour entries for each person (rows 3-8)
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 6

$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 6

$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 6

$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 6

$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 6

$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 6

# New TOTAL column header and per-row totals
$ws.Range("R1").Value = "TOTAL"
$ws.Range("R3").Formula = "=SUM(B3:Q3)"
$ws.Range("R4").Formula = "=SUM(B4:Q4)"
$ws.Range("R5").Formula = "=SUM(B5:Q5)"
$ws.Range("R6").Formula = "=SUM(B6:Q6)"
$ws.Range("R7").Formula = "=SUM(B7:Q7)"
$ws.Range("R8").Formula = "=SUM(B8:Q8)"

# Move the sheet selection to where the user last worked
$ws.Range("H9").Select()
